$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.614623665809631
$ws.Range("B1").Value = 1.799251675605774
$ws.Range("C1").Value = 1.82830798625946
$ws.Range("D1").Value = 2.371037721633911
$ws.Range("E1").Value = 3.852776050567627
